$wb = $excel.ActiveWorkbook

# --- Sheet R1 ---
$ws1 = $wb.Worksheets.Item("R1")
$ws1.Range("G2").Value = "3929:43:11"
$ws1.Range("G3").Value = "69:15:49"
$ws1.Range("D5").Value = "JED0125"
$ws1.Range("I5").Value = "Generator-SG"
$ws1.Range("J5").Value = "Good+In progress"

# --- Sheet R2 ---
$ws2 = $wb.Worksheets.Item("R2")
$ws2.Range("G2").Value = "12111:06:27"
$ws2.Range("G3").Value = "3240:49:56"
$ws2.Range("G4").Value = "479:01:30"

# --- Sheet R4 ---
$ws4 = $wb.Worksheets.Item("R4")
$ws4.Range("G2").Value = "2956:56:16"
$ws4.Range("G3").Value = "184:08:31"
$ws4.Range("G4").Value = "72:20:56"
$ws4.Range("G5").Value = "69:58:29"

# --- Sheet R5 ---
$ws5 = $wb.Worksheets.Item("R5")
$ws5.Range("G2").Value = "430:55:15"

# --- Sheet R6 ---
$ws6 = $wb.Worksheets.Item("R6")
$ws6.Range("G2").Value = "71:27:33"
